$d = $word.ActiveDocument

function Force-Split($startPos, $chunks) {
    # Given the absolute character offset of the start of a contiguous block of
    # text, and the ordered list of substrings that make it up, make sure each
    # substring ends up in its own run by toggling a character formatting
    # property on then back off over its exact range. (Toggling formatting
    # forces this engine to keep the edited sub-range as a distinct run instead
    # of folding it back into identically-formatted neighboring runs.)
    $pos = $startPos
    foreach ($c in $chunks) {
        $len = $c.Length
        if ($len -gt 0) {
            $cr = $d.Range($pos, $pos + $len)
            $cr.Font.Bold = 1
            $cr.Font.Bold = 0
        }
        $pos = $pos + $len
    }
}

# ---------------------------------------------------------------------------
# 1) Split "are" into "a" | bookmark(_GoBack) | "re" inside
#    "...decisive moves and positions are provided..."
#    (the bookmark naturally keeps the two pieces of "are" in separate runs)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("positions are provided")
$areStart = $rng.Start + "positions ".Length
$splitPoint = $areStart + 1

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 2) Replace every whole word "section" with "chapter" (same length: 7
#    characters so no other offsets shift), then restore original paragraph
#    run boundaries (this engine tends to merge all same-formatted runs in a
#    paragraph whenever any text inside it is edited).
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$positions = New-Object System.Collections.ArrayList
$searchFrom = 0
while ($true) {
    $idx = $full.IndexOf("section", $searchFrom)
    if ($idx -lt 0) { break }
    [void]$positions.Add($idx)
    $searchFrom = $idx + 7
}

Write-Host "Found $($positions.Count) occurrences of 'section'"

foreach ($pos in $positions) {
    $r = $d.Range($pos, $pos + 7)
    $r.Text = "chapter"
}

# ---------------------------------------------------------------------------
# 3) Re-establish the exact run layout for every paragraph touched above.
# ---------------------------------------------------------------------------

$full = $d.Content.Text

# -- Paragraph: "...decisive moves and positions a|re provided..." ----------
$p1Start = $full.IndexOf("    Chess games are noted")
$p1Chunks = @(
    "    ",
    "Chess games are noted in a uniform form like the PGN format in order to be able to perform analyses and evaluations at a later date.",
    " ",
    "These data usually contain information about the players and the event as well as the exact move sequence of the chess game.",
    " ",
    "In addition, especially in large chess databases, decisive moves and positions a",
    "re provided with comments after the game has been analyzed by grandmasters.",
    " ",
    "For a better arrangement of these comments there are standardized symbols and NAGs, which for example directly indicate whether a move was good or bad. These have the advantage that they are easier to evaluate due to their clear categorization and are also generally understandable, which is not the case with comments in natural language.",
    " \\\\"
)
Force-Split $p1Start $p1Chunks

# -- Paragraph: "The thesis is split into five more chapters: \\" -----------
$p2Start = $full.IndexOf("    The thesis is split into")
$p2Chunks = @(
    "    ",
    "The thesis is split into ",
    "five more ",
    "chapter",
    "s: \\\\"
)
Force-Split $p2Start $p2Chunks

# -- Paragraph: "The second chapter ... suitable ... such problems." --------
$p3Start = $full.IndexOf("    The second chapter")
$p3Chunks = @(
    "    ",
    "The second ",
    "chapter",
    " ",
    "initially provides background knowledge in the topics dealt with in this thesis.",
    " It starts with the ",
    "presentation of sentiment analysis as part of text mining",
    ". ",
    "Afterwards the concept of ",
    "w",
    "ord ",
    "e",
    "mbeddings is discussed, which can serve as a data model in such text mining problems.",
    " ",
    "This is followed by a description of multiclass classification problems and several ",
    "suitable",
    " approaches to solv",
    "e",
    " such problems. A separate look is taken at the subgroup of ordinal classification problems.  Finally, cost-sensitive methods are presented that take into account different weightings of misclassifications.",
    " \\\\"
)
Force-Split $p3Start $p3Chunks

# -- Paragraph: "The third chapter describes ... selected." -----------------
$p4Start = $full.IndexOf("    The third chapter")
$p4Chunks = @(
    "    ",
    "The third ",
    "chapter",
    " describes the general concept with which a text mining process can be created. It first deals with the requirements of problem and goal definition and the criteria for a suitable data selection. For the preparation of the data and their transformation into a model suitable for analysis, methods of natural language processing are presented.",
    " ",
    "Finally, possible analysis methods and evaluation methods are presented from which suitable techniques for the problem can be selected.",
    " \\\\"
)
Force-Split $p4Start $p4Chunks

# -- Paragraph: "In the fourth chapter, ... third chapter is applied ..." ---
$p5Start = $full.IndexOf("    In the fourth chapter")
$p5Chunks = @(
    "    ",
    "In the fourth ",
    "chapter",
    ", the procedure presented in the third ",
    "chapter",
    " is applied to the text mining process in chess annotations. First the format PGN and the annotation symbols NAG are explained and the five problems are specified. In the following ",
    "the",
    " used tools like NLTK, Weka and further libraries are listed. Finally, the classifiers and evaluation methods used are mentioned.",
    " \\\\"
)
Force-Split $p5Start $p5Chunks

# -- Paragraph: "The fifth chapter contains ... evaluation." ----------------
$p6Start = $full.IndexOf("    The fifth chapter")
$p6Chunks = @(
    "    ",
    "The fifth ",
    "chapter",
    " contains all results and their evaluations. After the analysis of different tokenizer configurations, statistics about the comments are generated, which can be used to gain basic knowledge about the analyzed data set. In addition, the attributes and models used for the data set are evaluated. The majority of the results are finally taken up by the comparisons of the achieved accuracies for all configurations and problems. In the end, the best results achieved are checked for optimality in a cost-sensitive evaluation.",
    " \\\\"
)
Force-Split $p6Start $p6Chunks

Write-Host "Done"
